{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is `async (context) => { ... }`.\n\nconst body = context.document.body;\n\n// Helper: find the single exact match of `oldText` in the document body\n// and replace it with `newText`.\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      `replaceOnce: expected exactly 1 match for ${JSON.stringify(oldText)}, found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// --- \"Relevant skills\" bullet list ---------------------------------------\n\nawait replaceOnce(\n  \"General: Object-oriented, but with a functional mindset. TDD/BDD.\",\n  \"General: Object-oriented, with a functional mindset. Focus on using TDD/BDD to drive clean, maintainable code.\"\n);\n\nawait replaceOnce(\n  \"Java: Spring, Camel, Java8 features, concurrency/threading, JMS, TDD/BDD (JUnit, Mockito, Cucumber-JVM), Servlets, REST (Jersey), Maven, JAXB, Jackson, Protobuf, Guava, Commons, Hibernate, DropWizard.\",\n  \"Java (5+ years): Spring, Camel, Java8 features, concurrency/threading, JMS, TDD/BDD (JUnit, Mockito, Cucumber-JVM), Servlets, REST (JAX-RS/Jersey), Maven, JAXB, Jackson, Protobuf, Guava, Commons, Hibernate, DropWizard.\"\n);\n\nawait replaceOnce(\n  \"Perl: Moose, Plack, DBIx::Class, Carton, XML, module development\",\n  'Perl (10+ years): \"Modern Perl\" (Moose/Moo, Plack etc), DBIx::Class, Carton, XML, Mason, mod_perl.'\n);\n\nawait replaceOnce(\n  \"Ruby: REST, XML, Sinatra, rvm/rbenv, gem development, bundler\",\n  \"Ruby (2 years for BDD): REST, XML, Sinatra, rvm/rbenv, gem development, bundler\"\n);\n\nawait replaceOnce(\n  \"Javascript: EC6, node.js\",\n  \"Javascript (occasional): EC6, node.js\"\n);\n\nawait replaceOnce(\n  \"BDD: Ruby, Cucumber, Jasmine\",\n  \"BDD: Ruby/Java, Cucumber, some Jasmine/Selenium\"\n);\n\n// --- BBC Future Media - Video Factory --------------------------------------\n\nawait replaceOnce(\n  \"Java/Camel/Jersey applications, continuous delivery, AWS, DevOps. Maintenance of legacy Perl stack.\",\n  \"Java/Camel/Jersey applications, continuous delivery, AWS, DevOps.\"\n);\n\n// Insert new bullet after \"Led a new team...\" (same bullet list).\n{\n  const results = body.search(\n    \"Led a new team, charged with reduction of both AWS costs and technical debt.\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"could not find 'Led a new team...' paragraph\");\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  para.insertParagraph(\"Maintenance of legacy Perl apps and CI.\", \"After\");\n  await context.sync();\n}\n\n// --- BBC Future Media - Publishing Services --------------------------------\n\nawait replaceOnce(\n  \"Development and maintenance of Java and perl components.\",\n  \"Development and maintenance of Java and Perl components.\"\n);\n\nawait replaceOnce(\n  'Design and delivery of \"Media Selector 5\" application in Perl; a business-critical service behind all iPlayer playback.',\n  'Design and delivery of \"Media Selector 5\" application in Perl; a business-critical, audience-facing service behind all iPlayer playback. Tech: Perl (Moose,Plack), XML/JSON, NoSQL (couchdb).'\n);\n\n// Insert new bullet after the (just updated) \"Media Selector 5\" bullet.\n{\n  const results = body.search(\n    'Design and delivery of \"Media Selector 5\" application in Perl',\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error('could not find \"Media Selector 5\" paragraph');\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  para.insertParagraph(\n    'Design and delivery of first iteration of \"Workflow Engine\" a business-critical publishing workflow for iPlayer content. Tech: Perl (Mason), Apache/mod_perl, HTTP, mySQL, DBIx::Class.',\n    \"After\"\n  );\n  await context.sync();\n}\n\n// --- CAIW Netwerken (contract) ---------------------------------------------\n\nawait replaceOnce(\n  \"Maintaining perl middleware services; reverse-engineered SOAP/WSDL contract from perl source code for Java clients.\",\n  \"Maintaining Perl middleware services; reverse-engineered SOAP/WSDL contract from Perl source code for Java clients.\"\n);\n\n// --- Semantico Ltd -----------------------------------------------------------\n\nawait replaceOnce(\n  \"Developing new client sites (e.g. www.blackwellreference.com) in perl using Semantico's backend CMS and access management technology.\",\n  \"Developing new client sites (e.g. www.blackwellreference.com) in Perl using Semantico's backend CMS and access management technology.\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($doc, [string]$oldText, [string]$newText) {\n    # Search the whole document body for an exact (case-sensitive), single\n    # occurrence of $oldText and overwrite it with $newText. We assign the\n    # found Range's .Text directly (rather than using Find.Execute's\n    # Replace:= argument) so straight quotes in $newText are not\n    # autocorrected into curly/smart quotes.\n    $rng = $doc.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop - do not wrap, so we can detect \"not found\"\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Replace-Once: text not found: $oldText\"\n    }\n    $rng.Text = $newText\n}\n\nfunction InsertBulletAfter($doc, [string]$anchorText, [string]$newText) {\n    # Find the paragraph containing $anchorText (an exact, single match) and\n    # insert a brand-new paragraph immediately after it, inheriting the\n    # anchor paragraph's style/list (pStyle + numId), with $newText as its\n    # content.\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $anchorText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"InsertBulletAfter: anchor text not found: $anchorText\"\n    }\n    $para = $find.Parent.Paragraphs(1)\n    $para.Range.InsertParagraphAfter()\n    $newPara = $para.Next()\n    $newPara.Range.Text = $newText\n}\n\n# --- \"Relevant skills\" bullet list ------------------------------------------\n\nReplace-Once $d `\n    \"General: Object-oriented, but with a functional mindset. TDD/BDD.\" `\n    \"General: Object-oriented, with a functional mindset. Focus on using TDD/BDD to drive clean, maintainable code.\"\n\nReplace-Once $d `\n    \"Java: Spring, Camel, Java8 features, concurrency/threading, JMS, TDD/BDD (JUnit, Mockito, Cucumber-JVM), Servlets, REST (Jersey), Maven, JAXB, Jackson, Protobuf, Guava, Commons, Hibernate, DropWizard.\" `\n    \"Java (5+ years): Spring, Camel, Java8 features, concurrency/threading, JMS, TDD/BDD (JUnit, Mockito, Cucumber-JVM), Servlets, REST (JAX-RS/Jersey), Maven, JAXB, Jackson, Protobuf, Guava, Commons, Hibernate, DropWizard.\"\n\nReplace-Once $d `\n    \"Perl: Moose, Plack, DBIx::Class, Carton, XML, module development\" `\n    'Perl (10+ years): \"Modern Perl\" (Moose/Moo, Plack etc), DBIx::Class, Carton, XML, Mason, mod_perl.'\n\nReplace-Once $d `\n    \"Ruby: REST, XML, Sinatra, rvm/rbenv, gem development, bundler\" `\n    \"Ruby (2 years for BDD): REST, XML, Sinatra, rvm/rbenv, gem development, bundler\"\n\nReplace-Once $d `\n    \"Javascript: EC6, node.js\" `\n    \"Javascript (occasional): EC6, node.js\"\n\nReplace-Once $d `\n    \"BDD: Ruby, Cucumber, Jasmine\" `\n    \"BDD: Ruby/Java, Cucumber, some Jasmine/Selenium\"\n\n# --- BBC Future Media - Video Factory ---------------------------------------\n\nReplace-Once $d `\n    \"Java/Camel/Jersey applications, continuous delivery, AWS, DevOps. Maintenance of legacy Perl stack.\" `\n    \"Java/Camel/Jersey applications, continuous delivery, AWS, DevOps.\"\n\nInsertBulletAfter $d `\n    \"Led a new team, charged with reduction of both AWS costs and technical debt.\" `\n    \"Maintenance of legacy Perl apps and CI.\"\n\n# --- BBC Future Media - Publishing Services ---------------------------------\n\nReplace-Once $d `\n    \"Development and maintenance of Java and perl components.\" `\n    \"Development and maintenance of Java and Perl components.\"\n\nReplace-Once $d `\n    'Design and delivery of \"Media Selector 5\" application in Perl; a business-critical service behind all iPlayer playback.' `\n    'Design and delivery of \"Media Selector 5\" application in Perl; a business-critical, audience-facing service behind all iPlayer playback. Tech: Perl (Moose,Plack), XML/JSON, NoSQL (couchdb).'\n\nInsertBulletAfter $d `\n    'Design and delivery of \"Media Selector 5\" application in Perl' `\n    'Design and delivery of first iteration of \"Workflow Engine\" a business-critical publishing workflow for iPlayer content. Tech: Perl (Mason), Apache/mod_perl, HTTP, mySQL, DBIx::Class.'\n\n# --- CAIW Netwerken (contract) -----------------------------------------------\n\nReplace-Once $d `\n    \"Maintaining perl middleware services; reverse-engineered SOAP/WSDL contract from perl source code for Java clients.\" `\n    \"Maintaining Perl middleware services; reverse-engineered SOAP/WSDL contract from Perl source code for Java clients.\"\n\n# --- Semantico Ltd ------------------------------------------------------------\n\nReplace-Once $d `\n    \"Developing new client sites (e.g. www.blackwellreference.com) in perl using Semantico's backend CMS and access management technology.\" `\n    \"Developing new client sites (e.g. www.blackwellreference.com) in Perl using Semantico's backend CMS and access management technology.\"\n"}
